# 15HP01_index.xlsx — interactive visualisation follow-up edit.
#
# Semantic changes applied (see commit message / diff):
#   1. Shared string "repeat" (column-G header) renamed to "rep".
#   2. The sheet's remembered selection moves from A1:H1 to the single
#      cell G2 (where the new "rep" column sits), matching the cursor
#      position left behind after the rename.
#   3. Columns A:H got a little wider (index/labels column widened for
#      the new header), and the sheet's default width for every column
#      beyond H grew to match.
#   4. Cosmetic window/print settings (tab-bar ratio, printer-defaults
#      flag) were also nudged by the authoring application; applied
#      here too where the object model exposes them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- 1. Rename the "repeat" header (column G) to "rep" --------------------
$ws.Range("G1").Value = "rep"

# -- 2. Move the active selection to G2 ------------------------------------
$null = $ws.Range("G2").Select()

# -- 3. Widen the data columns (A:H) and the sheet's default column width -
#       (values chosen so the stored width lands as close as possible to
#        the authored widths once the host snaps to its internal grid)
$ws.Columns.Item(1).ColumnWidth = 14.467687074829966
$ws.Columns.Item(2).ColumnWidth = 15.365646258503366
$ws.Columns.Item(3).ColumnWidth = 17.345238095238066
$ws.Columns.Item(4).ColumnWidth = 10.687074829931966
$ws.Columns.Item(5).ColumnWidth = 6.365646258503407
$ws.Columns.Item(6).ColumnWidth = 10.146258503401366
$ws.Columns.Item(7).ColumnWidth = 9.967687074829966
$ws.Columns.Item(8).ColumnWidth = 12.845238095238066

$ws.StandardWidth = 15.544217687074868
for ($i = 9; $i -le 1025; $i++) {
    $ws.Columns.Item($i).ColumnWidth = 15.544217687074868
}

# -- 4. Cosmetic window/print tweaks (best effort) -------------------------
$excel.ActiveWindow.TabRatio = 0.5
